# Add new "criteria grade" columns (Capex ... Waste and decomissioning) to
# the "Script (Main)" sheet, with a second header row of "/10" grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Script (Main)")

$headers = @(
    "Capex",
    "Safety",
    "Rentability",
    "Opex",
    "Ecological impact",
    "Startup time",
    "Scalability",
    "Availability (h/year)",
    "Plant Area/Footprint",
    "Technology readiness",
    "Connection flexibility",
    "Geopolitical barriers",
    "Economic lifetime",
    "Production efficiency",
    "Waste and decomissioning"
)

# Columns AF (32) through AT (46)
$startCol = 32
$endCol = $startCol + $headers.Count - 1

# Copy formatting from the last existing header/grade cells (AE1 / AE2) onto
# the new range so the new cells carry the same style indices.
$ws.Range("AE1").Copy()
$ws.Range($ws.Cells.Item(1, $startCol), $ws.Cells.Item(1, $endCol)).PasteSpecial(-4122)

$ws.Range("AE2").Copy()
$ws.Range($ws.Cells.Item(2, $startCol), $ws.Cells.Item(2, $endCol)).PasteSpecial(-4122)

# Write all the header labels (row 1) first so the shared-string table
# fills in header order, matching how the workbook was authored.
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Then the "/10" grade placeholders (row 2) so "/10" lands as the very last
# new shared string.
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = "/10"
}

# "Availability (h/year)" grade was actually entered as the number 10.
$ws.Cells.Item(2, 40).Value = 10

$excel.CutCopyMode = 0

# Widen the new columns the same way the author manually resized them.
# (ColumnWidth values chosen so the saved width lands as close as possible
# to the authored width.)
$ws.Columns.Item(34).ColumnWidth = 11.0
$ws.Columns.Item(36).ColumnWidth = 13.0
$ws.Columns.Item(39).ColumnWidth = 7.833333333333333
$ws.Columns.Item(40).ColumnWidth = 11.833333333333334
$ws.Columns.Item(41).ColumnWidth = 13.833333333333334
$ws.Columns.Item(42).ColumnWidth = 14.666666666666666
$ws.Columns.Item(43).ColumnWidth = 16.666666666666668
$ws.Columns.Item(44).ColumnWidth = 20.0
$ws.Columns.Item(45).ColumnWidth = 13.5
$ws.Columns.Item(46).ColumnWidth = 13.666666666666666

# Reflect the author's final view state: scrolled right to the new columns,
# with AN10 selected.
$ws.Activate()
$ws.Range("AN10").Select()
$excel.ActiveWindow.ScrollColumn = 30
